$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28 (leve id 27772)
$ws.Cells.Item(28, 8).Value = 434.35294
$ws.Cells.Item(28, 9).Value = 392.375
$ws.Cells.Item(28, 10).Value = 1106
$ws.Cells.Item(28, 11).Value = 392.375
$ws.Cells.Item(28, 12).Value = 1106
$ws.Cells.Item(28, 13).Value = 92.625
$ws.Cells.Item(28, 14).Value = -2076
# row 100 (leve id 19906)
$ws.Cells.Item(100, 8).Value = 1424.1111
$ws.Cells.Item(100, 9).Value = 1250.8334
$ws.Cells.Item(100, 10).Value = 1770.6666
$ws.Cells.Item(100, 11).Value = 1250.8334
$ws.Cells.Item(100, 12).Value = 1770.6666
$ws.Cells.Item(100, 13).Value = -709.8334
$ws.Cells.Item(100, 14).Value = -2852.6666
# row 132 (leve id 44049)
$ws.Cells.Item(132, 8).Value = 21564.318
$ws.Cells.Item(132, 9).Value = 3297.3242
$ws.Cells.Item(132, 10).Value = 89152.2
$ws.Cells.Item(132, 11).Value = 9891.972600000001
$ws.Cells.Item(132, 12).Value = 267456.6
$ws.Cells.Item(132, 13).Value = -7361.972600000001
$ws.Cells.Item(132, 14).Value = -272516.6

$ws = $wb.Worksheets.Item("ARM")
# row 7 (leve id 27125)
$ws.Cells.Item(7, 8).Value = 50000
$ws.Cells.Item(7, 10).Value = 50000
$ws.Cells.Item(7, 12).Value = 50000
$ws.Cells.Item(7, 14).Value = -50228
# row 19 (leve id 3550)
$ws.Cells.Item(19, 8).Value = 10100
$ws.Cells.Item(19, 9).Value = 0
$ws.Cells.Item(19, 10).Value = 10100
$ws.Cells.Item(19, 11).Value = 0
$ws.Cells.Item(19, 12).Value = 10100
$ws.Cells.Item(19, 13).ClearContents()
$ws.Cells.Item(19, 14).Value = -10558
# row 32 (leve id 44147)
$ws.Cells.Item(32, 8).Value = 29033.203
$ws.Cells.Item(32, 9).Value = 29673.309
$ws.Cells.Item(32, 10).Value = 24278.143
$ws.Cells.Item(32, 11).Value = 29673.309
$ws.Cells.Item(32, 12).Value = 24278.143
$ws.Cells.Item(32, 13).Value = -29386.309
$ws.Cells.Item(32, 14).Value = -24852.143
# row 61 (leve id 43999)
$ws.Cells.Item(61, 8).Value = 2764.2104
$ws.Cells.Item(61, 9).Value = 1626.4706
$ws.Cells.Item(61, 10).Value = 3685.238
$ws.Cells.Item(61, 11).Value = 1626.4706
$ws.Cells.Item(61, 12).Value = 3685.238
$ws.Cells.Item(61, 13).Value = -1414.4706
$ws.Cells.Item(61, 14).Value = -4109.237999999999
# row 102 (leve id 19945)
$ws.Cells.Item(102, 8).Value = 9183.134
$ws.Cells.Item(102, 9).Value = 1784.1177
$ws.Cells.Item(102, 10).Value = 18858.77
$ws.Cells.Item(102, 11).Value = 1784.1177
$ws.Cells.Item(102, 12).Value = 18858.77
$ws.Cells.Item(102, 13).Value = -162.1177
$ws.Cells.Item(102, 14).Value = -22102.77
# row 136 (leve id 43999)
$ws.Cells.Item(136, 8).Value = 2764.2104
$ws.Cells.Item(136, 9).Value = 1626.4706
$ws.Cells.Item(136, 10).Value = 3685.238
$ws.Cells.Item(136, 11).Value = 4879.4118
$ws.Cells.Item(136, 12).Value = 11055.714
$ws.Cells.Item(136, 13).Value = -2329.4118
$ws.Cells.Item(136, 14).Value = -16155.714

$ws = $wb.Worksheets.Item("BSM")
# row 19 (leve id 1753)
$ws.Cells.Item(19, 8).Value = 27648.334
$ws.Cells.Item(19, 9).Value = 1500
$ws.Cells.Item(19, 10).Value = 32878
$ws.Cells.Item(19, 11).Value = 1500
$ws.Cells.Item(19, 12).Value = 32878
$ws.Cells.Item(19, 13).Value = -1327
$ws.Cells.Item(19, 14).Value = -33224
# row 99 (leve id 19943)
$ws.Cells.Item(99, 8).Value = 2324.0667
$ws.Cells.Item(99, 9).Value = 2166.6667
$ws.Cells.Item(99, 10).Value = 2953.6667
$ws.Cells.Item(99, 11).Value = 2166.6667
$ws.Cells.Item(99, 12).Value = 2953.6667
$ws.Cells.Item(99, 13).Value = -668.6667000000002
$ws.Cells.Item(99, 14).Value = -5949.6667
# row 134 (leve id 43998)
$ws.Cells.Item(134, 8).Value = 3452.279
$ws.Cells.Item(134, 9).Value = 1508.36
$ws.Cells.Item(134, 11).Value = 4525.08
$ws.Cells.Item(134, 13).Value = -1990.08

$ws = $wb.Worksheets.Item("CRP")
# row 19 (leve id 2233)
$ws.Cells.Item(19, 8).Value = 489.66666
$ws.Cells.Item(19, 10).Value = 1300
$ws.Cells.Item(19, 12).Value = 1300
$ws.Cells.Item(19, 14).Value = -1640
# row 24 (leve id 2233)
$ws.Cells.Item(24, 8).Value = 489.66666
$ws.Cells.Item(24, 10).Value = 1300
$ws.Cells.Item(24, 12).Value = 1300
$ws.Cells.Item(24, 14).Value = -1640
# row 31 (leve id 44023)
$ws.Cells.Item(31, 8).Value = 186689.45
$ws.Cells.Item(31, 9).Value = 2085.9375
$ws.Cells.Item(31, 10).Value = 235110.05
$ws.Cells.Item(31, 11).Value = 2085.9375
$ws.Cells.Item(31, 12).Value = 235110.05
$ws.Cells.Item(31, 13).Value = -1790.9375
$ws.Cells.Item(31, 14).Value = -235700.05
# row 34 (leve id 44023)
$ws.Cells.Item(34, 8).Value = 186689.45
$ws.Cells.Item(34, 9).Value = 2085.9375
$ws.Cells.Item(34, 10).Value = 235110.05
$ws.Cells.Item(34, 11).Value = 2085.9375
$ws.Cells.Item(34, 12).Value = 235110.05
$ws.Cells.Item(34, 13).Value = -1883.9375
$ws.Cells.Item(34, 14).Value = -235514.05
# row 38 (leve id 1637)
$ws.Cells.Item(38, 8).Value = 32666.666
$ws.Cells.Item(38, 9).Value = 50000
$ws.Cells.Item(38, 10).Value = 24000
$ws.Cells.Item(38, 11).Value = 50000
$ws.Cells.Item(38, 12).Value = 24000
$ws.Cells.Item(38, 13).Value = -49623
$ws.Cells.Item(38, 14).Value = -24754
# row 46 (leve id 1637)
$ws.Cells.Item(46, 8).Value = 32666.666
$ws.Cells.Item(46, 9).Value = 50000
$ws.Cells.Item(46, 10).Value = 24000
$ws.Cells.Item(46, 11).Value = 50000
$ws.Cells.Item(46, 12).Value = 24000
$ws.Cells.Item(46, 13).Value = -49789
$ws.Cells.Item(46, 14).Value = -24422
# row 132 (leve id 44019)
$ws.Cells.Item(132, 8).Value = 45834.812
$ws.Cells.Item(132, 9).Value = 1738.3889
$ws.Cells.Item(132, 10).Value = 102530.21
$ws.Cells.Item(132, 11).Value = 5215.1667
$ws.Cells.Item(132, 12).Value = 307590.63
$ws.Cells.Item(132, 13).Value = -2685.1667
$ws.Cells.Item(132, 14).Value = -312650.63
# row 134 (leve id 44020)
$ws.Cells.Item(134, 8).Value = 53531.555
$ws.Cells.Item(134, 9).Value = 1150.6471
$ws.Cells.Item(134, 10).Value = 142579.1
$ws.Cells.Item(134, 11).Value = 3451.9413
$ws.Cells.Item(134, 12).Value = 427737.3
$ws.Cells.Item(134, 13).Value = -916.9412999999995
$ws.Cells.Item(134, 14).Value = -432807.3

$ws = $wb.Worksheets.Item("CUL")
# row 42 (leve id 4670)
$ws.Cells.Item(42, 8).Value = 702.5
$ws.Cells.Item(42, 10).Value = 702.5
$ws.Cells.Item(42, 12).Value = 2107.5
$ws.Cells.Item(42, 14).Value = -3175.5
# row 76 (leve id 12869)
$ws.Cells.Item(76, 8).Value = 4727.5474
$ws.Cells.Item(76, 9).Value = 2000
$ws.Cells.Item(76, 11).Value = 6000
$ws.Cells.Item(76, 13).Value = -5617
# row 79 (leve id 12869)
$ws.Cells.Item(79, 8).Value = 4727.5474
$ws.Cells.Item(79, 9).Value = 2000
$ws.Cells.Item(79, 11).Value = 6000
$ws.Cells.Item(79, 13).Value = -4674
# row 131 (leve id 36060)
$ws.Cells.Item(131, 8).Value = 3770
$ws.Cells.Item(131, 9).Value = 7593.2856
$ws.Cells.Item(131, 10).Value = 1711.3077
$ws.Cells.Item(131, 11).Value = 22779.8568
$ws.Cells.Item(131, 12).Value = 5133.9231
$ws.Cells.Item(131, 13).Value = -17739.8568
$ws.Cells.Item(131, 14).Value = -15213.9231

$ws = $wb.Worksheets.Item("GSM")
# row 6 (leve id 2108)
$ws.Cells.Item(6, 8).Value = 25747.5
$ws.Cells.Item(6, 9).Value = 25000
$ws.Cells.Item(6, 10).Value = 25996.666
$ws.Cells.Item(6, 11).Value = 25000
$ws.Cells.Item(6, 12).Value = 25996.666
$ws.Cells.Item(6, 13).Value = -24887
$ws.Cells.Item(6, 14).Value = -26222.666
# row 16 (leve id 2108)
$ws.Cells.Item(16, 8).Value = 25747.5
$ws.Cells.Item(16, 9).Value = 25000
$ws.Cells.Item(16, 10).Value = 25996.666
$ws.Cells.Item(16, 11).Value = 25000
$ws.Cells.Item(16, 12).Value = 25996.666
$ws.Cells.Item(16, 13).Value = -24750
$ws.Cells.Item(16, 14).Value = -26496.666
# row 40 (leve id 4113)
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(40, 10).Value = 0
$ws.Cells.Item(40, 12).Value = 0
$ws.Cells.Item(40, 14).ClearContents()
# row 46 (leve id 2078)
$ws.Cells.Item(46, 8).Value = 33692.25
$ws.Cells.Item(46, 10).Value = 33692.25
$ws.Cells.Item(46, 12).Value = 33692.25
$ws.Cells.Item(46, 14).Value = -34004.25
# row 70 (leve id 14146)
$ws.Cells.Item(70, 8).Value = 5222
$ws.Cells.Item(70, 9).Value = 5684.591
$ws.Cells.Item(70, 10).Value = 3949.875
$ws.Cells.Item(70, 11).Value = 5684.591
$ws.Cells.Item(70, 12).Value = 3949.875
$ws.Cells.Item(70, 13).Value = -5414.591
$ws.Cells.Item(70, 14).Value = -4489.875
# row 73 (leve id 14146)
$ws.Cells.Item(73, 8).Value = 5222
$ws.Cells.Item(73, 9).Value = 5684.591
$ws.Cells.Item(73, 10).Value = 3949.875
$ws.Cells.Item(73, 11).Value = 5684.591
$ws.Cells.Item(73, 12).Value = 3949.875
$ws.Cells.Item(73, 13).Value = -4748.591
$ws.Cells.Item(73, 14).Value = -5821.875
# row 132 (leve id 44008)
$ws.Cells.Item(132, 8).Value = 4097.2905
$ws.Cells.Item(132, 9).Value = 1532.625
$ws.Cells.Item(132, 10).Value = 6832.933
$ws.Cells.Item(132, 11).Value = 4597.875
$ws.Cells.Item(132, 12).Value = 20498.799
$ws.Cells.Item(132, 13).Value = -2067.875
$ws.Cells.Item(132, 14).Value = -25558.799

$ws = $wb.Worksheets.Item("LTW")
# row 119 (leve id 26288)
$ws.Cells.Item(119, 8).Value = 47408
$ws.Cells.Item(119, 10).Value = 47408
$ws.Cells.Item(119, 12).Value = 47408
$ws.Cells.Item(119, 14).Value = -57084
# row 132 (leve id 44058)
$ws.Cells.Item(132, 8).Value = 4136.7407
$ws.Cells.Item(132, 9).Value = 1914.2858
$ws.Cells.Item(132, 10).Value = 4914.6
$ws.Cells.Item(132, 11).Value = 5742.857400000001
$ws.Cells.Item(132, 12).Value = 14743.8
$ws.Cells.Item(132, 13).Value = -3212.857400000001
$ws.Cells.Item(132, 14).Value = -19803.8
# row 136 (leve id 44060)
$ws.Cells.Item(136, 8).Value = 1655.5312
$ws.Cells.Item(136, 9).Value = 1140.591
$ws.Cells.Item(136, 10).Value = 2788.4
$ws.Cells.Item(136, 11).Value = 3421.773
$ws.Cells.Item(136, 12).Value = 8365.200000000001
$ws.Cells.Item(136, 13).Value = -871.7729999999997
$ws.Cells.Item(136, 14).Value = -13465.2

$ws = $wb.Worksheets.Item("WVR")
# row 25 (leve id 3064)
$ws.Cells.Item(25, 8).Value = 20000
$ws.Cells.Item(25, 10).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 14).Value = -20586
# row 132 (leve id 44029)
$ws.Cells.Item(132, 8).Value = 1262.25
$ws.Cells.Item(132, 9).Value = 749.7692
$ws.Cells.Item(132, 10).Value = 3483
$ws.Cells.Item(132, 11).Value = 2249.3076
$ws.Cells.Item(132, 12).Value = 10449
$ws.Cells.Item(132, 13).Value = 280.6923999999999
$ws.Cells.Item(132, 14).Value = -15509
# row 136 (leve id 44031)
$ws.Cells.Item(136, 8).Value = 30238.39
$ws.Cells.Item(136, 9).Value = 65150.613
$ws.Cells.Item(136, 10).Value = 1757.3684
$ws.Cells.Item(136, 11).Value = 195451.839
$ws.Cells.Item(136, 12).Value = 5272.1052
$ws.Cells.Item(136, 13).Value = -192901.839
$ws.Cells.Item(136, 14).Value = -10372.1052
